$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17: evaluated cards count and "Aktiv" ranking value
$ws.Range("C17").Value = 11
$ws.Range("G17").Value = 10

# Row 20: percentage text update
$ws.Range("D20").Value = "90.9% der Karten"

# Rows 27-31: updated counts
$ws.Range("C27").Value = 105
$ws.Range("C28").Value = 105
$ws.Range("C29").Value = 105
$ws.Range("C30").Value = 102
$ws.Range("C31").Value = 16

# Row 34: percentage of cards without member
$ws.Range("G34").Value = "(63.6%)"

# Row 35 & 39: Theresa Schmid's counts
$ws.Range("C35").Value = 4
$ws.Range("F39").Value = 4

# Row 70 & 73: card totals
$ws.Range("B70").Value = 11
$ws.Range("B73").Value = 5
